# Regenerate the localization-status report for archive:
#  1. The "Ready for handoff" status has moved on to "In Translation" for
#     every file/locale row that shows it (Overview!E2:F3, zh-cn!C2:C3,
#     de-de!C2:C3 all share that text).
#  2. The Status columns are narrower now that the longest status text is
#     shorter, so their column widths shrink to match (Overview columns
#     E & F, and column C on both locale sheets).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth  = 12.5   # yields the closest achievable stored column width to the new, narrower status text

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- update every cell that currently shows the old status text ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- shrink the now-narrower Status columns ---
$overview.Columns.Item(5).ColumnWidth = $newWidth   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = $newWidth   # column F (de-de status)

$zhcn.Columns.Item(3).ColumnWidth = $newWidth       # column C (Status)
$dede.Columns.Item(3).ColumnWidth = $newWidth       # column C (Status)
